$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels for the new reporting date (20200816 -> 20200823)
$ws.Range("F1").Value = "AC_20200823"
$ws.Range("G1").Value = "AN_20200823"
$ws.Range("N1").Value = "AN_COTA_20200823"

# Update data values per the refreshed 26/08/2020 chart data extract
$ws.Range("C2").Value = -17.9
$ws.Range("D2").Value = -10.9
$ws.Range("E2").Value = 30.6
$ws.Range("F2").Value = 175.2
$ws.Range("G2").Value = 183.7
$ws.Range("N2").Value = -1.82
$ws.Range("C3").Value = -26.6
$ws.Range("D3").Value = -31.5
$ws.Range("E3").Value = 15.7
$ws.Range("F3").Value = 121.1
$ws.Range("G3").Value = 175.9
$ws.Range("N3").Value = -7.83
$ws.Range("C4").Value = -30.8
$ws.Range("D4").Value = -25.7
$ws.Range("E4").Value = 14.2
$ws.Range("F4").Value = 161.5
$ws.Range("G4").Value = 105
$ws.Range("C5").Value = -7.1
$ws.Range("D5").Value = 1.7
$ws.Range("E5").Value = 18.6
$ws.Range("F5").Value = 132.1
$ws.Range("G5").Value = 109.4
$ws.Range("C6").Value = -21.4
$ws.Range("D6").Value = -23.8
$ws.Range("E6").Value = -13.4
$ws.Range("F6").Value = 140.7
$ws.Range("G6").Value = 57.9
$ws.Range("N6").Value = -9.300000000000001
$ws.Range("C7").Value = -16.8
$ws.Range("D7").Value = -18.2
$ws.Range("E7").Value = 2.2
$ws.Range("F7").Value = 155.9
$ws.Range("G7").Value = 62
$ws.Range("C8").Value = -49.3
$ws.Range("D8").Value = -37.2
$ws.Range("E8").Value = 5.8
$ws.Range("F8").Value = 151.1
$ws.Range("G8").Value = 79.40000000000001
$ws.Range("C9").Value = -24.6
$ws.Range("D9").Value = -21.9
$ws.Range("E9").Value = -10.9
$ws.Range("F9").Value = 146.4
$ws.Range("G9").Value = 62.2
$ws.Range("C10").Value = -46.5
$ws.Range("D10").Value = -50.1
$ws.Range("E10").Value = -20.6
$ws.Range("F10").Value = 116.8
$ws.Range("G10").Value = 51
$ws.Range("N10").Value = 0.09
$ws.Range("C11").Value = -43.1
$ws.Range("D11").Value = -44.8
$ws.Range("E11").Value = -12.4
$ws.Range("F11").Value = 104.8
$ws.Range("G11").Value = 23.7
$ws.Range("N11").Value = 18.34
$ws.Range("C12").Value = -43.3
$ws.Range("D12").Value = -45.8
$ws.Range("E12").Value = -13
$ws.Range("F12").Value = 106.2
$ws.Range("G12").Value = 22.1
$ws.Range("N12").Value = 18.34
$ws.Range("C13").Value = -21.2
$ws.Range("D13").Value = -22.2
$ws.Range("E13").Value = -11.5
$ws.Range("F13").Value = 137.6
$ws.Range("G13").Value = 40.8
$ws.Range("C14").Value = -48.6
$ws.Range("D14").Value = -51.3
$ws.Range("E14").Value = -10.3
$ws.Range("F14").Value = 94
$ws.Range("G14").Value = 30
$ws.Range("N14").Value = 22
$ws.Range("C15").Value = -19.8
$ws.Range("D15").Value = -16.6
$ws.Range("E15").Value = 23.2
$ws.Range("F15").Value = 108.8
$ws.Range("G15").Value = 12.8
$ws.Range("N15").Value = 1.65
$ws.Range("C16").Value = -15.1
$ws.Range("D16").Value = -4
$ws.Range("E16").Value = 35.7
$ws.Range("F16").Value = 131.5
$ws.Range("G16").Value = 22.2
$ws.Range("N16").Value = 7.96
$ws.Range("C17").Value = -45.5
$ws.Range("D17").Value = -37.7
$ws.Range("E17").Value = 0.7
$ws.Range("F17").Value = 101.5
$ws.Range("G17").Value = 14
$ws.Range("C18").Value = -33.6
$ws.Range("D18").Value = -26.3
$ws.Range("E18").Value = 6.5
$ws.Range("F18").Value = 178.6
$ws.Range("G18").Value = 106.3
$ws.Range("N18").Value = -1.82
$ws.Range("C19").Value = -24.9
$ws.Range("D19").Value = -11
$ws.Range("E19").Value = 22.5
$ws.Range("F19").Value = 142.3
$ws.Range("G19").Value = 46.7
$ws.Range("C20").Value = -32.5
$ws.Range("D20").Value = -24.2
$ws.Range("E20").Value = 20.3
$ws.Range("F20").Value = 206.5
$ws.Range("G20").Value = 143.8
$ws.Range("C21").Value = -15.9
$ws.Range("D21").Value = -2.3
$ws.Range("E21").Value = 38.2
$ws.Range("F21").Value = 123.1
$ws.Range("G21").Value = 11
$ws.Range("C22").Value = -24.9
$ws.Range("D22").Value = -9.300000000000001
$ws.Range("E22").Value = 42.4
$ws.Range("F22").Value = 199.3
$ws.Range("G22").Value = 99.59999999999999
$ws.Range("C23").Value = -23.1
$ws.Range("D23").Value = -22.9
$ws.Range("E23").Value = -9.1
$ws.Range("F23").Value = 136.6
$ws.Range("G23").Value = 41.9
$ws.Range("C24").Value = -26.4
$ws.Range("D24").Value = -26.5
$ws.Range("E24").Value = -8.199999999999999
$ws.Range("F24").Value = 156.6
$ws.Range("G24").Value = 108.2
$ws.Range("C25").Value = -41.3
$ws.Range("D25").Value = -47.5
$ws.Range("E25").Value = -24.5
$ws.Range("F25").Value = 130.1
$ws.Range("G25").Value = 53.1
$ws.Range("C26").Value = -25
$ws.Range("D26").Value = -21.5
$ws.Range("E26").Value = 26
$ws.Range("F26").Value = 138.1
$ws.Range("G26").Value = 122.6
$ws.Range("C27").Value = -50
$ws.Range("D27").Value = -46
$ws.Range("E27").Value = -10.6
$ws.Range("F27").Value = 107
$ws.Range("G27").Value = 15.3
$ws.Range("C28").Value = -11.1
$ws.Range("D28").Value = -12.1
$ws.Range("E28").Value = 8.6
$ws.Range("F28").Value = 132.6
$ws.Range("G28").Value = 71.2
$ws.Range("C29").Value = -44.6
$ws.Range("D29").Value = -47.6
$ws.Range("E29").Value = -6.7
$ws.Range("F29").Value = 90.90000000000001
$ws.Range("G29").Value = -7.5
$ws.Range("N29").Value = 10.17
$ws.Range("C30").Value = -4.2
$ws.Range("D30").Value = 6.7
$ws.Range("E30").Value = 52.1
$ws.Range("F30").Value = 145.2
$ws.Range("G30").Value = 42.7
$ws.Range("C31").Value = -14.8
$ws.Range("D31").Value = -14.9
$ws.Range("E31").Value = 1.1
$ws.Range("F31").Value = 133.6
$ws.Range("G31").Value = 56.3
$ws.Range("C32").Value = -10.5
$ws.Range("D32").Value = -11.6
$ws.Range("E32").Value = 6.2
$ws.Range("F32").Value = 132.3
$ws.Range("G32").Value = 72
$ws.Range("C33").Value = -30.6
$ws.Range("D33").Value = -27.7
$ws.Range("E33").Value = -8.5
$ws.Range("F33").Value = 96.8
$ws.Range("G33").Value = 12.1
$ws.Range("C34").Value = -48
$ws.Range("D34").Value = -50.3
$ws.Range("E34").Value = -8.4
$ws.Range("F34").Value = 91.59999999999999
$ws.Range("G34").Value = 32.4
$ws.Range("N34").Value = 22
$ws.Range("C35").Value = -52.1
$ws.Range("D35").Value = -47.3
$ws.Range("E35").Value = -12.1
$ws.Range("F35").Value = 111.7
$ws.Range("G35").Value = 10.7
$ws.Range("C36").Value = -57.5
$ws.Range("D36").Value = -56.6
$ws.Range("E36").Value = -18.5
$ws.Range("F36").Value = 71.3
$ws.Range("G36").Value = -31.5
$ws.Range("M36").Value = -78.22
$ws.Range("N36").Value = -37.54
$ws.Range("C37").Value = -8.300000000000001
$ws.Range("D37").Value = -4.9
$ws.Range("E37").Value = 30.7
$ws.Range("F37").Value = 130.7
$ws.Range("G37").Value = 42.3
$ws.Range("C38").Value = -38.3
$ws.Range("D38").Value = -43.2
$ws.Range("E38").Value = -8.699999999999999
$ws.Range("F38").Value = 129.8
$ws.Range("G38").Value = 64.40000000000001
$ws.Range("N38").Value = 0.09
$ws.Range("C39").Value = -2.6
$ws.Range("D39").Value = 10.6
$ws.Range("E39").Value = 56.7
$ws.Range("F39").Value = 146.6
$ws.Range("G39").Value = 37
$ws.Range("C40").Value = -10.5
$ws.Range("D40").Value = -11.1
$ws.Range("E40").Value = 18.7
$ws.Range("F40").Value = 129.5
$ws.Range("G40").Value = 51.5
$ws.Range("C41").Value = -46.5
$ws.Range("D41").Value = -46
$ws.Range("E41").Value = -10.9
$ws.Range("F41").Value = 102.2
$ws.Range("G41").Value = 26.9
$ws.Range("N41").Value = 22
$ws.Range("C42").Value = -23.5
$ws.Range("D42").Value = -28.7
$ws.Range("E42").Value = -16.7
$ws.Range("F42").Value = 156.6
$ws.Range("G42").Value = 60.6
$ws.Range("N42").Value = -9.300000000000001
